{"js": "// Under \"CH\u01af\u01a0NG 2: C\u01a0 S\u1ede L\u00dd THUY\u1ebeT\", three ListParagraph sub-heading items\n// (numId=8) are promoted to bold, outline-level-2 sub-headings. The third\n// one (currently empty) also gets its missing heading text filled in.\n\nconst body = context.document.body;\n\n// Locate the three target paragraphs by their (unique) text so the script\n// does not depend on brittle absolute paragraph indices.\nconst firstResults = body.search(\"GI\u1edaI THI\u1ec6U V\u1ec0 STM32F103C8T6 V\u00c0 FOTA\", { matchCase: true });\nconst secondResults = body.search(\"GI\u1edaI THI\u1ec6U V\u1ec0 GIAO TH\u1ee8C HTTP\", { matchCase: true });\nfirstResults.load(\"items\");\nsecondResults.load(\"items\");\nawait context.sync();\n\nconst firstPara = firstResults.items[0].paragraphs.getFirst();\nconst secondPara = secondResults.items[0].paragraphs.getFirst();\nconst thirdPara = secondPara.getNext();\n\n// The third paragraph is empty in the original document; give it its\n// missing heading text before formatting it like its two siblings.\nthirdPara.insertText(\n  \"PH\u00c2N CHIA B\u1ed8 NH\u1eda CH\u01af\u01a0NG TR\u00ccNH V\u00c0 CH\u01af\u01a0NG TR\u00ccNH BOOTLOADER\",\n  Word.InsertLocation.replace\n);\n\nfor (const para of [firstPara, secondPara, thirdPara]) {\n  // OOXML <w:outlineLvl w:val=\"1\"/> <-> Word OM OutlineLevel 2.\n  para.outlineLevel = 2;\n  // Bold the whole paragraph (runs + paragraph mark), matching <w:b/>/<w:bCs/>.\n  para.font.bold = true;\n  para.font.boldBidirectional = true;\n}\n\nawait context.sync();\n", "ps1": "# Under \"CH\u01af\u01a0NG 2: C\u01a0 S\u1ede L\u00dd THUY\u1ebeT\", three ListParagraph sub-heading items\n# (numId=8) are promoted to bold, outline-level-2 sub-headings. The third\n# one (currently empty) also gets its missing heading text filled in.\n\n$d = $word.ActiveDocument\n\n# Locate the first two target paragraphs by their (unique) text so the\n# script does not depend on brittle absolute paragraph indices.\n$findRange1 = $d.Content\n$findRange1.Find.Execute(\"GI\u1edaI THI\u1ec6U V\u1ec0 STM32F103C8T6 V\u00c0 FOTA\") | Out-Null\n$firstPara = $findRange1.Paragraphs(1)\n\n$findRange2 = $d.Content\n$findRange2.Find.Execute(\"GI\u1edaI THI\u1ec6U V\u1ec0 GIAO TH\u1ee8C HTTP\") | Out-Null\n$secondPara = $findRange2.Paragraphs(1)\n\n$thirdPara = $secondPara.Next()\n\n# The third paragraph is empty in the original document; give it its\n# missing heading text before formatting it like its two siblings.\n$thirdPara.Range.Text = \"PH\u00c2N CHIA B\u1ed8 NH\u1eda CH\u01af\u01a0NG TR\u00ccNH V\u00c0 CH\u01af\u01a0NG TR\u00ccNH BOOTLOADER\"\n\nforeach ($para in @($firstPara, $secondPara, $thirdPara)) {\n    # OOXML <w:outlineLvl w:val=\"1\"/> <-> Word OM OutlineLevel 2.\n    $para.OutlineLevel = 2\n    # Bold the whole paragraph (runs + paragraph mark), matching <w:b/>/<w:bCs/>.\n    $para.Range.Font.Bold = 1\n    $para.Range.Font.BoldBi = 1\n}\n"}
